$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-16 18:32:41"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
